$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: widen column L ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns.Item(12).ColumnWidth = 5.6666666666667

# --- Transactions sheet: update values and selection ---
$wsTx = $wb.Worksheets.Item("Transactions")

$wsTx.Range("A2").Value = 61
$wsTx.Range("J2").Value = 9133.2199999999993

$wsTx.Range("A3").Value = 59
$wsTx.Range("C3").Value = 42064
$wsTx.Range("E3").Value = 963.77
$wsTx.Range("F3").Value = 866.78
$wsTx.Range("G3").Value = 96.99
$wsTx.Range("J3").Value = 4133.22

$wsTx.Range("A4").Value = 57

$wsTx.Activate()
$wsTx.Range("A2:L4").Select()
